# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated numbers (gh-pages output refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (all exhibition rows)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 44
$ws1.Range("F3").Value = 321
$ws1.Range("F4").Value = 248
$ws1.Range("F5").Value = 2952
$ws1.Range("F6").Value = 2002
$ws1.Range("F7").Value = 386
$ws1.Range("F9").Value = 1093
$ws1.Range("F11").Value = 547
$ws1.Range("F12").Value = 58

# Sheet "全部类型" (combined listing, includes one extra 演出 row at row 8)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 44
$ws4.Range("F3").Value = 321
$ws4.Range("F4").Value = 248
$ws4.Range("F5").Value = 2952
$ws4.Range("F6").Value = 2002
$ws4.Range("F7").Value = 386
$ws4.Range("F10").Value = 1093
$ws4.Range("F12").Value = 547
$ws4.Range("F13").Value = 58
